$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Developer" column (E) values were stored in lower-case ("michael",
# "jason", "jake"). Re-capitalise them to "Michael", "Jason", "Jake" while
# leaving every other cell (including the Notes column) untouched.
for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -eq "michael") {
        $cell.Value = "Michael"
    } elseif ($val -eq "jason") {
        $cell.Value = "Jason"
    } elseif ($val -eq "jake") {
        $cell.Value = "Jake"
    }
}

# Update the saved selection so it covers the whole table (A1:E12) instead
# of the previous single-cell selection (H7).
$ws.Range("A1:E12").Select()
